# Insert a new data row at row 223 (pushing existing rows 223..272 down to
# 224..273) and populate it with the new weekly price observation, matching
# the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 223; Excel copies formatting
# (including the date number format on column D) from the row above by
# default, which matches the style="2" used on the rest of column D.
$ws.Rows.Item(223).Insert()

$ws.Cells.Item(223, 1).Value = 7
$ws.Cells.Item(223, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(223, 3).Value = "Ñuble"
$ws.Cells.Item(223, 4).Value = 44932
$ws.Cells.Item(223, 5).Value = 16
$ws.Cells.Item(223, 6).Value = 100112032
$ws.Cells.Item(223, 7).Value = "Zapallo italiano"
$ws.Cells.Item(223, 8).Value = "Sin especificar"
$ws.Cells.Item(223, 9).Value = "Primera"
$ws.Cells.Item(223, 10).Value = 120
$ws.Cells.Item(223, 11).Value = 5500
$ws.Cells.Item(223, 12).Value = 6000
$ws.Cells.Item(223, 13).Value = 5750
$ws.Cells.Item(223, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(223, 15).Value = "Región del Maule"
$ws.Cells.Item(223, 16).Value = 115
$ws.Cells.Item(223, 17).Value = 50
$ws.Cells.Item(223, 18).Value = "Hortaliza"
